$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 115
$ws.Range("G6").Value = 3436.2
$ws.Range("F9").Value = 15
$ws.Range("G9").Value = 443.55
$ws.Range("B10").Value = 33532.43
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("F77").Value = 299
$ws.Range("G77").Value = 13975.26
$ws.Range("F79").Value = 75
$ws.Range("G79").Value = 4662.75
$ws.Range("F83").Value = 122
$ws.Range("G83").Value = 18381.74
$ws.Range("B90").Value = 209557.76
$ws.Range("B127").Value = 57552
$ws.Range("E127").Value = 136.86
$ws.Range("F127").Value = -5
$ws.Range("G127").Value = -603.45
$ws.Range("B128").Value = 64329
$ws.Range("E128").Value = 128.32
$ws.Range("F128").Value = 2
$ws.Range("G128").Value = 241.38
$ws.Range("F144").Value = 1600
$ws.Range("G144").Value = 13520
$ws.Range("F146").Value = 37
$ws.Range("G146").Value = 3115.03
$ws.Range("B147").Value = 23378.59
$ws.Range("F150").Value = 48
$ws.Range("G150").Value = 2231.52
$ws.Range("B156").Value = 37196.27
$ws.Range("F167").Value = 22
$ws.Range("G167").Value = 6314.22
$ws.Range("B175").Value = 35268.88
$ws.Range("F203").Value = 80
$ws.Range("G203").Value = 1612.8
$ws.Range("F211").Value = 56
$ws.Range("G211").Value = 5667.2
$ws.Range("B216").Value = 56955.24
$ws.Range("F218").Value = 8
$ws.Range("G218").Value = 1729.76
$ws.Range("B219").Value = 61610
$ws.Range("E219").Value = 122.71
$ws.Range("F219").Value = -58
$ws.Range("G219").Value = -5957.18
$ws.Range("B220").Value = 63565
$ws.Range("E220").Value = 109.19
$ws.Range("F220").Value = 60
$ws.Range("G220").Value = 6162.6
$ws.Range("F223").Value = 18
$ws.Range("G223").Value = 2384.64
$ws.Range("F227").Value = 67
$ws.Range("G227").Value = 9666.76
$ws.Range("B232").Value = 63510
$ws.Range("E232").Value = 50.66
$ws.Range("F232").Value = 127
$ws.Range("G232").Value = 6050.28
$ws.Range("B233").Value = 55356
$ws.Range("E233").Value = 54.04
$ws.Range("F233").Value = -158
$ws.Range("G233").Value = -7527.12
$ws.Range("F240").Value = 6
$ws.Range("G240").Value = 3129.84
$ws.Range("F248").Value = 2
$ws.Range("G248").Value = 118.26
$ws.Range("F256").Value = 302
$ws.Range("G256").Value = 45653.34
$ws.Range("B260").Value = 216956.25
$ws.Range("F288").Value = 54
$ws.Range("G288").Value = 5021.46
$ws.Range("F291").Value = 131
$ws.Range("G291").Value = 5634.31
$ws.Range("F295").Value = 7
$ws.Range("G295").Value = 725.83
$ws.Range("B304").Value = 203745.61
$ws.Range("F306").Value = 75
$ws.Range("G306").Value = 1583.25
$ws.Range("B309").Value = 2005.02
$ws.Range("B322").Value = 47097
$ws.Range("D322").Value = 112.28
$ws.Range("E322").Value = 134.16
$ws.Range("F322").Value = 15
$ws.Range("G322").Value = 1684.2
$ws.Range("B323").Value = 58047
$ws.Range("D323").Value = 105.54
$ws.Range("E323").Value = 126.1
$ws.Range("F323").Value = 41
$ws.Range("G323").Value = 4327.14
$ws.Range("F328").Value = 69
$ws.Range("G328").Value = 2567.49
$ws.Range("B330").Value = 33033.37
$ws.Range("F334").Value = 201
$ws.Range("G334").Value = 10415.82
$ws.Range("F338").Value = 84
$ws.Range("G338").Value = 1990.8
$ws.Range("F339").Value = 10
$ws.Range("G339").Value = 474
$ws.Range("F343").Value = 42
$ws.Range("G343").Value = 3022.74
$ws.Range("F345").Value = 88
$ws.Range("G345").Value = 5404.08
$ws.Range("B346").Value = 30314.91
$ws.Range("B364").Value = 53602
$ws.Range("E364").Value = 15.69
$ws.Range("F364").Value = -231
$ws.Range("G364").Value = -3037.65
$ws.Range("B365").Value = 65068
$ws.Range("E365").Value = 13.97
$ws.Range("F365").Value = 63
$ws.Range("G365").Value = 828.45
$ws.Range("B380").Value = 64925
$ws.Range("E380").Value = 13.97
$ws.Range("F380").Value = 111
$ws.Range("G380").Value = 1459.65
$ws.Range("B381").Value = 45709
$ws.Range("E381").Value = 15.69
$ws.Range("F381").Value = -300
$ws.Range("G381").Value = -3945
$ws.Range("F456").Value = 51
$ws.Range("G456").Value = 5639.07
$ws.Range("B460").Value = 16779.52
$ws.Range("B463").Value = 60025
$ws.Range("E463").Value = 37.22
$ws.Range("F463").Value = -98
$ws.Range("G463").Value = -3217.34
$ws.Range("B464").Value = 64833
$ws.Range("E464").Value = 34.9
$ws.Range("F464").Value = 95
$ws.Range("G464").Value = 3118.85
$ws.Range("B473").Value = 60022
$ws.Range("E473").Value = 37.22
$ws.Range("F473").Value = -113
$ws.Range("G473").Value = -3709.79
$ws.Range("B474").Value = 64830
$ws.Range("E474").Value = 34.9
$ws.Range("F474").Value = 109
$ws.Range("G474").Value = 3578.47
$ws.Range("F477").Value = 18
$ws.Range("G477").Value = 816.12
$ws.Range("B478").Value = 816.12
$ws.Range("F485").Value = 30
$ws.Range("G485").Value = 5264.1
$ws.Range("F486").Value = 77
$ws.Range("G486").Value = 6802.18
$ws.Range("B488").Value = 33408.98
$ws.Range("F509").Value = 261
$ws.Range("G509").Value = 20979.18
$ws.Range("B510").Value = 28151.04
$ws.Range("F551").Value = 19
$ws.Range("G551").Value = 2719.47
$ws.Range("F554").Value = 18
$ws.Range("G554").Value = 671.04
$ws.Range("F558").Value = 44
$ws.Range("G558").Value = 5940.44
$ws.Range("B560").Value = 17522.95
$ws.Range("B572").Value = 65079
$ws.Range("F572").Value = 18
$ws.Range("G572").Value = 735.66
$ws.Range("B573").Value = 65362
$ws.Range("F573").Value = 29
$ws.Range("G573").Value = 1185.23
$ws.Range("F582").Value = 60
$ws.Range("G582").Value = 3419.4
$ws.Range("B583").Value = 32341.29
$ws.Range("F599").Value = 2278
$ws.Range("G599").Value = 371564.58
$ws.Range("F601").Value = 493
$ws.Range("G601").Value = 139454.91
$ws.Range("B606").Value = 565388.04
$ws.Range("F610").Value = 17
$ws.Range("G610").Value = 696.83
$ws.Range("B618").Value = 49378.72
$ws.Range("B619").Value = 2108023.03
$ws.Range("B620").Value = 2108023.03